$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "TI" (Trading Instrument) worksheet after the existing
#    "Logger" sheet and make it the active tab.
# ---------------------------------------------------------------------------
$logger = $wb.Worksheets.Item("Logger")
$ti = $wb.Worksheets.Add($null, $logger)
$ti.Name = "TI"

# ---------------------------------------------------------------------------
# 2. Populate the header / title rows.
#    (Cell write order below intentionally mirrors the order the strings
#    appear in the authored sharedStrings table.)
# ---------------------------------------------------------------------------
$ti.Range("A1:E1").Merge() | Out-Null
$ti.Range("A1").Value = "Description"
$ti.Range("A1:E1").ShrinkToFit = $false

$ti.Range("F1").Value = "Trading Instrument Requirements"

$ti.Range("A2").Value = "Status"
$ti.Range("B2").Value = "Priority"
$ti.Range("D2").Value = "Req#"
$ti.Range("E2").Value = "Category"
$ti.Range("F2").Value = "Description"

# ---------------------------------------------------------------------------
# 3. Requirement rows.
# ---------------------------------------------------------------------------
$ti.Range("E3").Value = "FUNC"
$ti.Range("F3").Value = "A list of available trading instruments shall be created from the app.config configuration file at startup. Invalid entries shall be ignored."

$ti.Range("E4").Value = "FUNC"
$ti.Range("F4").Value = "Available trading instruments shall automatically subscribe to real-time data after connection to the broker."
$ti.Range("F4").WrapText = $true

$ti.Range("C2").Value = "User Story"

$ti.Range("E5").Value = "FUNC"
$ti.Range("F5").Value = "Trading data shall be stored in 1-minute increments for the last year for each available instrument. The High/Low/Open/Close for each interval shall be stored in a separate file."

$ti.Range("E6").Value = "FUNC"
$ti.Range("F6").Value = 'Trading data shall be stored the following format: "TIMESTAMP,OPEN,CLOSE,HIGH,LOW". The timestamp shall be "MONTH/DAY/YEAR" with Month/Day as one or 2 digits and year as 4 digits.'

$ti.Range("E7").Value = "FUNC"
$ti.Range("F7").Value = 'Trading prices shall add zeros such that all data has the same string length (i.e: "1000.25,1005,1000.50,900,500.005" shall be coded as "1000.250,1005.000,1000.500,500.005"'

# ---------------------------------------------------------------------------
# 4. Column widths.
# ---------------------------------------------------------------------------
$ti.Columns.Item(1).ColumnWidth = 7.8776041666667
$ti.Columns.Item(2).ColumnWidth = 9.0221354166667
$ti.Columns.Item(3).ColumnWidth = 11.5924479166667
$ti.Columns.Item(4).ColumnWidth = 6.8776041666667
$ti.Columns.Item(5).ColumnWidth = 10.1666666666667
$ti.Columns.Item(6).ColumnWidth = 211.0221354166667

# ---------------------------------------------------------------------------
# 5. Data validation lists.
# ---------------------------------------------------------------------------
$ti.Range("E3:E49").Validation.Add(3, 1, 1, '"FUNC,NON-FUNC"') | Out-Null
$ti.Range("A3:C49").Validation.Add(3, 1, 1, '"Open,Done"') | Out-Null

# ---------------------------------------------------------------------------
# 6. Page setup.
# ---------------------------------------------------------------------------
$ti.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 7. Turn the header/data range into an Excel Table (ListObject).
# ---------------------------------------------------------------------------
$tiTable = $ti.ListObjects.Add(1, $ti.Range("A2:F49"), $null, 1)
$tiTable.Name = "Table13"
$tiTable.TableStyle = "TableStyleMedium13"

# ---------------------------------------------------------------------------
# 8. Activate the new sheet and select F10 to match the authored selection.
# ---------------------------------------------------------------------------
$ti.Activate()
$ti.Range("F10").Select() | Out-Null
